$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the existing header comments (row 15) before touching the
#     sheet. Columns 1..74 == A..BV. A few header cells (AL15, AM15, AU15)
#     intentionally have no comment, so we leave them absent from the map.
$origCount = 74
$origComments = @{}
for ($col = 1; $col -le $origCount; $col++) {
    $cell = $ws.Cells.Item(15, $col)
    $cmt = $cell.Comment
    if ($cmt -ne $null) {
        $origComments[$col] = $cmt.Text()
    }
}

# --- Insert a new column at U (21st column) - "culture_collection" is
#     being added between "chem_administration" (T) and "depth" (which
#     slides from U to V).
$ws.Columns("U").Insert()

# --- New header text/value for the inserted column.
$ws.Range("U15").Value2 = "culture_collection"

$newCommentText = "Name of source institute and unique culture identifier. See the description for the proper format and list of allowed institutes, http://www.insdc.org/controlled-vocabulary-culturecollection-qualifier"

# --- Shift the captured comments one column to the right (U..BV -> V..BW),
#     walking from the far end so we never overwrite a source we still need.
for ($col = $origCount; $col -ge 21; $col--) {
    $targetCol = $col + 1
    $text = $origComments[$col]
    $targetCell = $ws.Cells.Item(15, $targetCol)

    if ($text -ne $null) {
        if ($targetCell.Comment -ne $null) {
            [void]$targetCell.Comment.Text($text)
        } else {
            [void]$targetCell.AddComment($text)
        }
    } else {
        if ($targetCell.Comment -ne $null) {
            [void]$targetCell.Comment.Delete()
        }
    }
}


# --- Give the newly inserted column (U15) its own comment.
$uCell = $ws.Range("U15")
if ($uCell.Comment -ne $null) {
    [void]$uCell.Comment.Text($newCommentText)
} else {
    [void]$uCell.AddComment($newCommentText)
}
